$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 29600
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 29600
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 29600
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = -30186

$ws.Range("H98").Value = 2290.125
$ws.Range("I98").Value = 2187.8572
$ws.Range("K98").Value = 2187.8572
$ws.Range("M98").Value = -689.8571999999999

$ws.Range("H111").Value = 703.7857
$ws.Range("I111").Value = 711.7143
$ws.Range("J111").Value = 695.8570999999999
$ws.Range("K111").Value = 2135.1429
$ws.Range("L111").Value = 2087.5713
$ws.Range("M111").Value = 931.8571000000002
$ws.Range("N111").Value = -8221.5713

$ws.Range("H116").Value = 12241.833
$ws.Range("I116").Value = 15010
$ws.Range("J116").Value = 3937.3333
$ws.Range("K116").Value = 15010
$ws.Range("L116").Value = 3937.3333
$ws.Range("M116").Value = -11568
$ws.Range("N116").Value = -10821.3333

$ws.Range("H122").Value = 2290.125
$ws.Range("I122").Value = 2187.8572
$ws.Range("K122").Value = 6563.571599999999
$ws.Range("M122").Value = -4113.571599999999

$ws.Range("H135").Value = 2529.739
$ws.Range("J135").Value = 4337.5
$ws.Range("L135").Value = 39037.5
$ws.Range("N135").Value = -44107.5

$ws.Range("H138").Value = 2465.761
$ws.Range("I138").Value = 1399.2812
$ws.Range("J138").Value = 3034.55
$ws.Range("K138").Value = 4197.8436
$ws.Range("L138").Value = 9103.650000000001
$ws.Range("M138").Value = 942.1563999999998
$ws.Range("N138").Value = -19383.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 85006.25
$ws.Range("J23").Value = 100007
$ws.Range("L23").Value = 100007
$ws.Range("N23").Value = -100525

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null

$ws.Range("H51").Value = 30000
$ws.Range("J51").Value = 30000
$ws.Range("L51").Value = 30000
$ws.Range("N51").Value = -31512

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = $null

$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("K57").Value = 5000
$ws.Range("M57").Value = -4516

$ws.Range("H122").Value = 734940.9
$ws.Range("I122").Value = 1027841.06
$ws.Range("J122").Value = 2690.5
$ws.Range("K122").Value = 3083523.18
$ws.Range("L122").Value = 8071.5
$ws.Range("M122").Value = -3081073.18
$ws.Range("N122").Value = -12971.5

$ws.Range("H132").Value = 1541106.1
$ws.Range("I132").Value = 1816.8
$ws.Range("J132").Value = 5004507
$ws.Range("K132").Value = 5450.4
$ws.Range("L132").Value = 15013521
$ws.Range("M132").Value = -2920.4
$ws.Range("N132").Value = -15018581

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 68333.336
$ws.Range("J122").Value = 68333.336
$ws.Range("L122").Value = 68333.336
$ws.Range("N122").Value = -78133.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 32251.25
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 32251.25
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 32251.25
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = -32529.25

$ws.Range("H14").Value = 10000
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10340

$ws.Range("H20").Value = 59800
$ws.Range("J20").Value = 59800
$ws.Range("L20").Value = 59800
$ws.Range("N20").Value = -60272

$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = -950

$ws.Range("H30").Value = 59800
$ws.Range("J30").Value = 59800
$ws.Range("L30").Value = 59800
$ws.Range("N30").Value = -59982

$ws.Range("H31").Value = 234265.31
$ws.Range("I31").Value = 1951
$ws.Range("J31").Value = 340321.84
$ws.Range("K31").Value = 1951
$ws.Range("L31").Value = 340321.84
$ws.Range("M31").Value = -1656
$ws.Range("N31").Value = -340911.84

$ws.Range("H34").Value = 234265.31
$ws.Range("I34").Value = 1951
$ws.Range("J34").Value = 340321.84
$ws.Range("K34").Value = 1951
$ws.Range("L34").Value = 340321.84
$ws.Range("M34").Value = -1749
$ws.Range("N34").Value = -340725.84

$ws.Range("H59").Value = 25883.334
$ws.Range("I59").Value = 30000
$ws.Range("J59").Value = 25060
$ws.Range("K59").Value = 30000
$ws.Range("L59").Value = 25060
$ws.Range("M59").Value = -28855
$ws.Range("N59").Value = -27350

$ws.Range("H128").Value = 59800
$ws.Range("J128").Value = 59800
$ws.Range("L128").Value = 59800
$ws.Range("N128").Value = -69760

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 233.4
$ws.Range("J15").Value = 241.75
$ws.Range("L15").Value = 725.25
$ws.Range("N15").Value = -1005.25

$ws.Range("H17").Value = 693.5
$ws.Range("I17").Value = 693.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2080.5
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -1911.5
$ws.Range("N17").Value = $null

$ws.Range("H46").Value = 780
$ws.Range("I46").Value = 300
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -809
$ws.Range("N46").Value = -4682

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5250
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -384
$ws.Range("N3").Value = -10232

$ws.Range("H5").Value = 8500
$ws.Range("J5").Value = 8500
$ws.Range("L5").Value = 8500
$ws.Range("N5").Value = -8724

$ws.Range("H10").Value = 32650000
$ws.Range("I10").Value = 48500000
$ws.Range("K10").Value = 48500000
$ws.Range("M10").Value = -48499831

$ws.Range("H11").Value = 24428572
$ws.Range("I11").Value = 25200000
$ws.Range("K11").Value = 25200000
$ws.Range("M11").Value = -25199861

$ws.Range("H97").Value = 2240.6155
$ws.Range("I97").Value = 2181.9
$ws.Range("J97").Value = 2436.3333
$ws.Range("K97").Value = 2181.9
$ws.Range("L97").Value = 2436.3333
$ws.Range("M97").Value = -1685.9
$ws.Range("N97").Value = -3428.3333

$ws.Range("H122").Value = 39440676
$ws.Range("I122").Value = 70989496
$ws.Range("K122").Value = 212968488
$ws.Range("M122").Value = -212966038

$ws.Range("H132").Value = 3221.1833
$ws.Range("I132").Value = 2968.4324
$ws.Range("J132").Value = 3627.7827
$ws.Range("K132").Value = 8905.297200000001
$ws.Range("L132").Value = 10883.3481
$ws.Range("M132").Value = -6375.297200000001
$ws.Range("N132").Value = -15943.3481

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 5749.75
$ws.Range("J11").Value = 5749.75
$ws.Range("L11").Value = 5749.75
$ws.Range("N11").Value = -6029.75

$ws.Range("H68").Value = 2075.3
$ws.Range("I68").Value = 1764.2858
$ws.Range("J68").Value = 2801
$ws.Range("K68").Value = 1764.2858
$ws.Range("L68").Value = 2801
$ws.Range("M68").Value = -1015.2858
$ws.Range("N68").Value = -4299

$ws.Range("H71").Value = 2075.3
$ws.Range("I71").Value = 1764.2858
$ws.Range("J71").Value = 2801
$ws.Range("K71").Value = 8821.429
$ws.Range("L71").Value = 14005
$ws.Range("M71").Value = -5077.429
$ws.Range("N71").Value = -21493

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 23750.5
$ws.Range("J2").Value = 23750.5
$ws.Range("L2").Value = 23750.5
$ws.Range("N2").Value = -23974.5

$ws.Range("H18").Value = 1000000000
$ws.Range("J18").Value = 1000000000
$ws.Range("L18").Value = 1000000000
$ws.Range("N18").Value = -1000000346

$ws.Range("H20").Value = 5003605
$ws.Range("J20").Value = 7210
$ws.Range("L20").Value = 7210
$ws.Range("N20").Value = -7690

$ws.Range("H81").Value = 1650.125
$ws.Range("I81").Value = 1166.8334
$ws.Range("K81").Value = 2333.6668
$ws.Range("M81").Value = -1272.6668

$ws.Range("H84").Value = 1650.125
$ws.Range("I84").Value = 1166.8334
$ws.Range("K84").Value = 11668.334
$ws.Range("M84").Value = -6364.333999999999

$ws.Range("H94").Value = 31332.5
$ws.Range("J94").Value = 31332.5
$ws.Range("L94").Value = 31332.5
$ws.Range("N94").Value = -33134.5

$ws.Range("H135").Value = 62557.5
$ws.Range("J135").Value = 62557.5
$ws.Range("L135").Value = 62557.5
$ws.Range("N135").Value = -72697.5

$ws.Range("H136").Value = 2310.4922
$ws.Range("I136").Value = 2198.054
$ws.Range("J136").Value = 2459.0715
$ws.Range("K136").Value = 6594.162
$ws.Range("L136").Value = 7377.2145
$ws.Range("M136").Value = -4044.162
$ws.Range("N136").Value = -12477.2145

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

$ws.Range("H139").Value = 46558.145
$ws.Range("J139").Value = 46558.145
$ws.Range("L139").Value = 46558.145
$ws.Range("N139").Value = -56838.145
